$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C8").Value = "عقوبة"
$ws.Range("C9").Value = "عقوبه"

$ws.Columns.Item(2).ColumnWidth = 21

$ws.Range("I27").Select()
